# Apply edits to REZISTIVITE (sheet1) and SISMIK (sheet2)
$wb = $excel.ActiveWorkbook

# --- REZISTIVITE sheet ---
$ws1 = $wb.Worksheets.Item("REZISTIVITE")

# For each data row (2..7), shift the values in columns E, G, I, K, M one
# step to the left (E<-G, G<-I, I<-K, K<-M, M<-empty).
$cols = @("E", "G", "I", "K", "M")
for ($row = 2; $row -le 7; $row++) {
    $vals = @()
    foreach ($col in $cols) {
        $vals += $ws1.Range("$col$row").Value2
    }
    for ($i = 0; $i -lt ($cols.Length - 1); $i++) {
        $ws1.Range("$($cols[$i])$row").Value = $vals[$i + 1]
    }
    $ws1.Range("$($cols[$cols.Length - 1])$row").Value = $null
}

# --- SISMIK sheet ---
$ws2 = $wb.Worksheets.Item("SISMIK")

$cols2 = @("E", "H", "K")
for ($row = 2; $row -le 9; $row++) {
    $vals = @()
    foreach ($col in $cols2) {
        $vals += $ws2.Range("$col$row").Value2
    }
    for ($i = 0; $i -lt ($cols2.Length - 1); $i++) {
        $ws2.Range("$($cols2[$i])$row").Value = $vals[$i + 1]
    }
    $ws2.Range("$($cols2[$cols2.Length - 1])$row").Value = $null
}

# Set per-sheet selections. Select SISMIK's range first, then REZISTIVITE's
# last so REZISTIVITE (originally tabSelected) ends up the active sheet
# again, matching the source workbook's tab selection state.
$ws2.Range("K2:K9").Select() | Out-Null
$ws1.Range("M9").Select() | Out-Null
